# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 2
# of the zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 01:02:20"
$wsZhCn.Range("H2").Value = "2016-03-22 01:02:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 01:02:24"
$wsDeDe.Range("H2").Value = "2016-03-22 01:03:01"
